# Update the EMTR supplemental table for the February 2023 baseline.
# Shifts the year header row and refreshes the data rows with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - year headers (2022-2032 -> 2023-2033)
$ws.Cells.Item(9, 2).Value = 2023
$ws.Cells.Item(9, 3).Value = 2024
$ws.Cells.Item(9, 4).Value = 2025
$ws.Cells.Item(9, 5).Value = 2026
$ws.Cells.Item(9, 6).Value = 2027
$ws.Cells.Item(9, 7).Value = 2028
$ws.Cells.Item(9, 8).Value = 2029
$ws.Cells.Item(9, 9).Value = 2030
$ws.Cells.Item(9, 10).Value = 2031
$ws.Cells.Item(9, 11).Value = 2032
$ws.Cells.Item(9, 12).Value = 2033

# Row 12 - Overall
$ws.Cells.Item(12, 2).Value = 15.28
$ws.Cells.Item(12, 3).Value = 15.63
$ws.Cells.Item(12, 4).Value = 15.74
$ws.Cells.Item(12, 5).Value = 15.62
$ws.Cells.Item(12, 6).Value = 16.44
$ws.Cells.Item(12, 7).Value = 16.46
$ws.Cells.Item(12, 8).Value = 16.24
$ws.Cells.Item(12, 9).Value = 16.26
$ws.Cells.Item(12, 10).Value = 16.23
$ws.Cells.Item(12, 11).Value = 16.32
$ws.Cells.Item(12, 12).Value = 16.49

# Row 14 - Business Assets
$ws.Cells.Item(14, 2).Value = 18.84
$ws.Cells.Item(14, 3).Value = 19.47
$ws.Cells.Item(14, 4).Value = 19.9
$ws.Cells.Item(14, 5).Value = 22.55
$ws.Cells.Item(14, 6).Value = 23.6
$ws.Cells.Item(14, 7).Value = 23.66
$ws.Cells.Item(14, 8).Value = 23.52
$ws.Cells.Item(14, 9).Value = 23.56
$ws.Cells.Item(14, 10).Value = 23.58
$ws.Cells.Item(14, 11).Value = 23.71
$ws.Cells.Item(14, 12).Value = 23.9

# Row 16 - Nonresidential equipment
$ws.Cells.Item(16, 2).Value = 9.84
$ws.Cells.Item(16, 3).Value = 12.12
$ws.Cells.Item(16, 4).Value = 14.27
$ws.Cells.Item(16, 5).Value = 17.2
$ws.Cells.Item(16, 6).Value = 20.28
$ws.Cells.Item(16, 7).Value = 20.35
$ws.Cells.Item(16, 8).Value = 20.21
$ws.Cells.Item(16, 9).Value = 20.25
$ws.Cells.Item(16, 10).Value = 20.27
$ws.Cells.Item(16, 11).Value = 20.38
$ws.Cells.Item(16, 12).Value = 20.57

# Row 17 - Nonresidential structures
$ws.Cells.Item(17, 2).Value = 21.02
$ws.Cells.Item(17, 3).Value = 21.27
$ws.Cells.Item(17, 4).Value = 21.27
$ws.Cells.Item(17, 5).Value = 23.22
$ws.Cells.Item(17, 6).Value = 23.88
$ws.Cells.Item(17, 7).Value = 23.94
$ws.Cells.Item(17, 8).Value = 23.81
$ws.Cells.Item(17, 9).Value = 23.85
$ws.Cells.Item(17, 10).Value = 23.87
$ws.Cells.Item(17, 11).Value = 23.99
$ws.Cells.Item(17, 12).Value = 24.18

# Row 18 - Residential property
$ws.Cells.Item(18, 2).Value = 25.29
$ws.Cells.Item(18, 3).Value = 25.38
$ws.Cells.Item(18, 4).Value = 25.09
$ws.Cells.Item(18, 5).Value = 29.2
$ws.Cells.Item(18, 6).Value = 29.39
$ws.Cells.Item(18, 7).Value = 29.41
$ws.Cells.Item(18, 8).Value = 29.25
$ws.Cells.Item(18, 9).Value = 29.31
$ws.Cells.Item(18, 10).Value = 29.31
$ws.Cells.Item(18, 11).Value = 29.46
$ws.Cells.Item(18, 12).Value = 29.64

# Row 19 - R&D and own-account software
$ws.Cells.Item(19, 2).Value = -4.93
$ws.Cells.Item(19, 3).Value = -5.16
$ws.Cells.Item(19, 4).Value = -5.39
$ws.Cells.Item(19, 5).Value = -4.14
$ws.Cells.Item(19, 6).Value = -4.09
$ws.Cells.Item(19, 7).Value = -3.98
$ws.Cells.Item(19, 8).Value = -4.14
$ws.Cells.Item(19, 9).Value = -4.06
$ws.Cells.Item(19, 10).Value = -3.98
$ws.Cells.Item(19, 11).Value = -3.78
$ws.Cells.Item(19, 12).Value = -3.5

# Row 20 - Other intellectual property products
$ws.Cells.Item(20, 2).Value = 14.52
$ws.Cells.Item(20, 3).Value = 18.13
$ws.Cells.Item(20, 4).Value = 21.45
$ws.Cells.Item(20, 5).Value = 25.17
$ws.Cells.Item(20, 6).Value = 28.21
$ws.Cells.Item(20, 7).Value = 28.28
$ws.Cells.Item(20, 8).Value = 28.15
$ws.Cells.Item(20, 9).Value = 28.18
$ws.Cells.Item(20, 10).Value = 28.2
$ws.Cells.Item(20, 11).Value = 28.31
$ws.Cells.Item(20, 12).Value = 28.5

# Row 21 - Inventories
$ws.Cells.Item(21, 2).Value = 30.6
$ws.Cells.Item(21, 3).Value = 30.54
$ws.Cells.Item(21, 4).Value = 30.33
$ws.Cells.Item(21, 5).Value = 33.15
$ws.Cells.Item(21, 6).Value = 33.26
$ws.Cells.Item(21, 7).Value = 33.31
$ws.Cells.Item(21, 8).Value = 33.2
$ws.Cells.Item(21, 9).Value = 33.24
$ws.Cells.Item(21, 10).Value = 33.26
$ws.Cells.Item(21, 11).Value = 33.39
$ws.Cells.Item(21, 12).Value = 33.55

# Row 23 - Equity-financed
$ws.Cells.Item(23, 2).Value = 20.94
$ws.Cells.Item(23, 3).Value = 21.57
$ws.Cells.Item(23, 4).Value = 22.07
$ws.Cells.Item(23, 5).Value = 25.16
$ws.Cells.Item(23, 6).Value = 26.15
$ws.Cells.Item(23, 7).Value = 26.21
$ws.Cells.Item(23, 8).Value = 26.14
$ws.Cells.Item(23, 9).Value = 26.15
$ws.Cells.Item(23, 10).Value = 26.17
$ws.Cells.Item(23, 11).Value = 26.27
$ws.Cells.Item(23, 12).Value = 26.37

# Row 24 - Debt-financed
$ws.Cells.Item(24, 2).Value = 9.59
$ws.Cells.Item(24, 3).Value = 10.26
$ws.Cells.Item(24, 4).Value = 10.34
$ws.Cells.Item(24, 5).Value = 10.29
$ws.Cells.Item(24, 6).Value = 11.57
$ws.Cells.Item(24, 7).Value = 11.63
$ws.Cells.Item(24, 8).Value = 11.19
$ws.Cells.Item(24, 9).Value = 11.41
$ws.Cells.Item(24, 10).Value = 11.44
$ws.Cells.Item(24, 11).Value = 11.74
$ws.Cells.Item(24, 12).Value = 12.37

# Row 26 - C corporations
$ws.Cells.Item(26, 2).Value = 17.58
$ws.Cells.Item(26, 3).Value = 18.2
$ws.Cells.Item(26, 4).Value = 18.76
$ws.Cells.Item(26, 5).Value = 19.88
$ws.Cells.Item(26, 6).Value = 21.04
$ws.Cells.Item(26, 7).Value = 21.13
$ws.Cells.Item(26, 8).Value = 21
$ws.Cells.Item(26, 9).Value = 21.03
$ws.Cells.Item(26, 10).Value = 21.06
$ws.Cells.Item(26, 11).Value = 21.17
$ws.Cells.Item(26, 12).Value = 21.36

# Row 27 - Pass-through entities
$ws.Cells.Item(27, 2).Value = 21.66
$ws.Cells.Item(27, 3).Value = 22.33
$ws.Cells.Item(27, 4).Value = 22.48
$ws.Cells.Item(27, 5).Value = 28.2
$ws.Cells.Item(27, 6).Value = 29.04
$ws.Cells.Item(27, 7).Value = 29.04
$ws.Cells.Item(27, 8).Value = 28.88
$ws.Cells.Item(27, 9).Value = 28.95
$ws.Cells.Item(27, 10).Value = 28.95
$ws.Cells.Item(27, 11).Value = 29.11
$ws.Cells.Item(27, 12).Value = 29.29

# Row 29 - Owner-Occupied Housing Structures
$ws.Cells.Item(29, 2).Value = 7.98
$ws.Cells.Item(29, 3).Value = 7.65
$ws.Cells.Item(29, 4).Value = 6.99
$ws.Cells.Item(29, 5).Value = -0.71
$ws.Cells.Item(29, 6).Value = -0.68
$ws.Cells.Item(29, 7).Value = -0.77
$ws.Cells.Item(29, 8).Value = -1.21
$ws.Cells.Item(29, 9).Value = -1.26
$ws.Cells.Item(29, 10).Value = -1.43
$ws.Cells.Item(29, 11).Value = -1.46
$ws.Cells.Item(29, 12).Value = -1.39

# Row 31 - Overall, Including Land (B and E:L only; C31/D31 unchanged)
$ws.Cells.Item(31, 2).Value = -0.23
$ws.Cells.Item(31, 5).Value = -3.32
$ws.Cells.Item(31, 6).Value = -3.38
$ws.Cells.Item(31, 7).Value = -3.38
$ws.Cells.Item(31, 8).Value = -3.37
$ws.Cells.Item(31, 9).Value = -3.37
$ws.Cells.Item(31, 10).Value = -3.37
$ws.Cells.Item(31, 11).Value = -3.37
$ws.Cells.Item(31, 12).Value = -3.37

# Row 32 - Business Assets, Including Land
$ws.Cells.Item(32, 2).Value = 26.85
$ws.Cells.Item(32, 3).Value = 26.01
$ws.Cells.Item(32, 4).Value = 24.12
$ws.Cells.Item(32, 5).Value = 6.87
$ws.Cells.Item(32, 6).Value = 7.19
$ws.Cells.Item(32, 7).Value = 6.85
$ws.Cells.Item(32, 8).Value = 5.15
$ws.Cells.Item(32, 9).Value = 4.93
$ws.Cells.Item(32, 10).Value = 4.28
$ws.Cells.Item(32, 11).Value = 4.18
$ws.Cells.Item(32, 12).Value = 4.5

# Row 35 - Overall (Owner-Occupied Housing, Including Land section)
$ws.Cells.Item(35, 2).Value = 15.8
$ws.Cells.Item(35, 3).Value = 16.02
$ws.Cells.Item(35, 4).Value = 15.97
$ws.Cells.Item(35, 5).Value = 15.49
$ws.Cells.Item(35, 6).Value = 16.13
$ws.Cells.Item(35, 7).Value = 16.14
$ws.Cells.Item(35, 8).Value = 15.9
$ws.Cells.Item(35, 9).Value = 15.92
$ws.Cells.Item(35, 10).Value = 15.88
$ws.Cells.Item(35, 11).Value = 15.97
$ws.Cells.Item(35, 12).Value = 16.13

# Row 36 - Equity-financed
$ws.Cells.Item(36, 2).Value = 20.65
$ws.Cells.Item(36, 3).Value = 21.14
$ws.Cells.Item(36, 4).Value = 21.42
$ws.Cells.Item(36, 5).Value = 24.23
$ws.Cells.Item(36, 6).Value = 25.09
$ws.Cells.Item(36, 7).Value = 25.14
$ws.Cells.Item(36, 8).Value = 25
$ws.Cells.Item(36, 9).Value = 25.05
$ws.Cells.Item(36, 10).Value = 25.06
$ws.Cells.Item(36, 11).Value = 25.19
$ws.Cells.Item(36, 12).Value = 25.38

# Row 37 - Debt-financed
$ws.Cells.Item(37, 2).Value = 27.97
$ws.Cells.Item(37, 3).Value = 27.97
$ws.Cells.Item(37, 4).Value = 27.69
$ws.Cells.Item(37, 5).Value = 31.07
$ws.Cells.Item(37, 6).Value = 31.2
$ws.Cells.Item(37, 7).Value = 31.23
$ws.Cells.Item(37, 8).Value = 31.09
$ws.Cells.Item(37, 9).Value = 31.15
$ws.Cells.Item(37, 10).Value = 31.16
$ws.Cells.Item(37, 11).Value = 31.3
$ws.Cells.Item(37, 12).Value = 31.48

# Row 38 - Business land
$ws.Cells.Item(38, 2).Value = 7.98
$ws.Cells.Item(38, 3).Value = 7.65
$ws.Cells.Item(38, 4).Value = 6.99
$ws.Cells.Item(38, 5).Value = -0.71
$ws.Cells.Item(38, 6).Value = -0.68
$ws.Cells.Item(38, 7).Value = -0.77
$ws.Cells.Item(38, 8).Value = -1.21
$ws.Cells.Item(38, 9).Value = -1.26
$ws.Cells.Item(38, 10).Value = -1.43
$ws.Cells.Item(38, 11).Value = -1.46
$ws.Cells.Item(38, 12).Value = -1.39
